# Generate Report for Handback
#
# The localization-status report is regenerated once a handback is fully
# in sync with en-US:
#   - the "Status" column (and the Overview roll-up columns that mirror
#     it) flips from "Ready for handoff" to "Handed back: in sync with
#     en-US" for every language sheet;
#   - each language sheet's "Latest Handback DateTime" is refreshed to
#     the timestamp of this handback;
#   - the stale "handback file is not the latest" error message is
#     cleared now that the handback is current;
#   - a couple of report columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Excel's Range.ColumnWidth setter only stores whole sixths of a
# character (it rounds to the nearest pixel internally), so asking for
# the exact fractional target width would overshoot to the wrong
# bucket. Back-solve for the input that lands in the bucket nearest the
# desired stored width.
function Set-PreciseColumnWidth($column, [double]$targetWidth) {
    $bucket = [Math]::Round($targetWidth * 6)
    $column.ColumnWidth = (($bucket - 5.5) / 6) + 0.08
}

# --- Overview sheet: handback status shown per language -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
Set-PreciseColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-PreciseColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527

# --- zh-cn handback sheet --------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-09-03 19:02:31"
$wsZhCn.Range("P2").Value = ""
Set-PreciseColumnWidth $wsZhCn.Columns.Item(3) 29.9777047293527
Set-PreciseColumnWidth $wsZhCn.Columns.Item(16) 13.7470528738839

# --- de-de handback sheet --------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-09-03 19:02:38"
$wsDeDe.Range("P2").Value = ""
Set-PreciseColumnWidth $wsDeDe.Columns.Item(3) 29.9777047293527
Set-PreciseColumnWidth $wsDeDe.Columns.Item(16) 13.7470528738839
